$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.07007951928338062
$ws.Range("C2").Value = 0.7205358658215948
$ws.Range("D2").Value = 0.7411812092429714
$ws.Range("E2").Value = 0.8609188168712375
$ws.Range("F2").Value = 0.88293831854245
$ws.Range("G2").Value = 18
